$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")
$ws.Columns.Item(1).ColumnWidth = 13.42578125
